# Adds two new checklist items (each: checkbox picture + " " + text) to the
# "Dhamma Sineru Deep Clean: Residences Checklist" document:
#   1. "Open all windows"   -> end of the Bedroom cell in the "2. Dry Dust" table
#   2. "Mattress stains"    -> end of the (a) Bedroom cell in the "3. Scrub" table
#
# Each existing checklist item paragraph is: <w:drawing> (unchecked checkbox,
# r:embed="rId23") + a run with a literal space + a run with the item text.
# We replicate that exact pattern via Range.InsertXML so the new paragraphs
# are indistinguishable in structure from their siblings.

function Get-ParaIndexContaining($cellRange, $substr) {
    $count = $cellRange.Paragraphs.Count
    $idx = -1
    for ($i = 1; $i -le $count; $i++) {
        $pp = $cellRange.Paragraphs($i)
        if ($pp.Range.Text -like "*$substr*") {
            $idx = $i
        }
    }
    return $idx
}

function Add-ChecklistItemAfter($cellRange, $anchorSubstr, $itemText) {
    $idx = Get-ParaIndexContaining $cellRange $anchorSubstr
    if ($idx -lt 0) {
        throw "Could not find anchor paragraph containing '$anchorSubstr'"
    }

    $anchorPara = $cellRange.Paragraphs($idx)
    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $cellRange.Paragraphs($idx + 1)
    $target = $newPara.Range
    $target.Collapse(1)

    $xmlFrag = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
 xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"
 xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"
 xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"
 xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
  <w:pPr><w:jc w:val="left"/></w:pPr>
  <w:r>
    <w:drawing>
      <wp:inline>
        <wp:extent cx="179043" cy="179043"/>
        <wp:effectExtent b="0" l="0" r="0" t="0"/>
        <wp:docPr descr="" title="" id="1" name="Picture"/>
        <a:graphic>
          <a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">
            <pic:pic>
              <pic:nvPicPr>
                <pic:cNvPr descr="images/unchecked.png" id="0" name="Picture"/>
                <pic:cNvPicPr>
                  <a:picLocks noChangeArrowheads="1" noChangeAspect="1"/>
                </pic:cNvPicPr>
              </pic:nvPicPr>
              <pic:blipFill>
                <a:blip r:embed="rId23"/>
                <a:stretch>
                  <a:fillRect/>
                </a:stretch>
              </pic:blipFill>
              <pic:spPr bwMode="auto">
                <a:xfrm>
                  <a:off x="0" y="0"/>
                  <a:ext cx="179043" cy="179043"/>
                </a:xfrm>
                <a:prstGeom prst="rect">
                  <a:avLst/>
                </a:prstGeom>
                <a:noFill/>
                <a:ln w="9525">
                  <a:noFill/>
                  <a:headEnd/>
                  <a:tailEnd/>
                </a:ln>
              </pic:spPr>
            </pic:pic>
          </a:graphicData>
        </a:graphic>
      </wp:inline>
    </w:drawing>
  </w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t xml:space="preserve">$itemText</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

    $target.InsertXML($xmlFrag)
}

$d = $word.ActiveDocument

# "2. Dry Dust" table -> Bedroom cell -> add "Open all windows" after
# "Bedframe (top and bottom)" (the last item in that cell).
$dryDustBedroomCell = $d.Tables(1).Cell(1, 1)
Add-ChecklistItemAfter $dryDustBedroomCell.Range "Bedframe (top and bottom)" "Open all windows"

# "3. Scrub" table -> (a) Bedroom cell -> add "Mattress stains" after
# "Shelves" (the last item in that cell).
$scrubBedroomCell = $d.Tables(2).Cell(1, 1)
Add-ChecklistItemAfter $scrubBedroomCell.Range "Shelves" "Mattress stains"

Write-Output "Done."
